$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing numeric values (columns B..F) for rows 2-9
$ws.Range("B2").Value = 0.2716253482700288
$ws.Range("C2").Value = 0.646407500340176
$ws.Range("D2").Value = 0.84144457876875
$ws.Range("E2").Value = 0.917302882786678
$ws.Range("F2").Value = 0.9092388822629066

$ws.Range("B3").Value = 0.2726843341995326
$ws.Range("C3").Value = 0.731747839592782
$ws.Range("D3").Value = 1.006886897929402
$ws.Range("E3").Value = 1.003437540621937
$ws.Range("F3").Value = 1.005107455845322

$ws.Range("B4").Value = 0.3732018941855291
$ws.Range("C4").Value = 0.72231839705317
$ws.Range("D4").Value = 0.8136630649843791
$ws.Range("E4").Value = 0.9020327405279585
$ws.Range("F4").Value = 0.8577243111825058

$ws.Range("B5").Value = 0.5036323863627127
$ws.Range("C5").Value = 0.6286073955106629
$ws.Range("D5").Value = 0.6937851054938032
$ws.Range("E5").Value = 0.832937636017129
$ws.Range("F5").Value = 0.6958113806129091

$ws.Range("B6").Value = 0.4544332582864191
$ws.Range("C6").Value = 0.5646968050727786
$ws.Range("D6").Value = 0.4534038351377512
$ws.Range("E6").Value = 0.6733526825800512
$ws.Range("F6").Value = 0.5237622964888429

$ws.Range("B7").Value = 0.3601308786133244
$ws.Range("C7").Value = 0.423148671095416
$ws.Range("D7").Value = 0.5203223525804069
$ws.Range("E7").Value = 0.7213337317638812
$ws.Range("F7").Value = 0.6629152402123534
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = 0.5515821818765249
$ws.Range("C8").Value = 0.600657948696442
$ws.Range("D8").Value = 0.6415195273476259
$ws.Range("E8").Value = 0.8009491415487164
$ws.Range("F8").Value = 0.6361854672819474
$ws.Range("G8").Value = 6

$ws.Range("B9").Value = 0.8660249170169981
$ws.Range("C9").Value = 0.8660249170169981
$ws.Range("D9").Value = 0.8404044078926498
$ws.Range("E9").Value = 0.9167357350363571
$ws.Range("F9").Value = 0.3682497474507306
$ws.Range("G9").Value = 3

# New row 10
$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = -0.04443783748577212
$ws.Range("C10").Value = 0.04443783748577212
$ws.Range("D10").Value = 0.001974721400411894
$ws.Range("E10").Value = 0.04443783748577212
$ws.Range("G10").Value = 1

# Copy the style from A9 (existing labeled cell) onto the new A10 label cell
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
